# Updates cryptos list values (Price / Volume(1h) columns, and one full row
# swap from "EnergySwap" to "Algorand") to match the latest scraped data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a plain-text value into a cell without letting Excel
# auto-convert numeric-looking strings (e.g. "0.992") into real numbers,
# and without leaving a lasting "Text" number-format style on the cell.
function Set-TextValue($range, [string]$text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

$ws.Range('D2').Value = '26.063.75'
$ws.Range('E2').Value = '  +1.09%  '
$ws.Range('D3').Value = '1.635.81'
$ws.Range('E3').Value = '  +0.06%  '
Set-TextValue $ws.Range('D4') '0.992'
$ws.Range('E4').Value = '  -0.93%  '
Set-TextValue $ws.Range('D5') '215.34'
$ws.Range('E5').Value = '  +0.09%  '
$ws.Range('E6').Value = '  -0.34%  '
$ws.Range('E7').Value = '  -0.88%  '
$ws.Range('E8').Value = '  -0.98%  '
Set-TextValue $ws.Range('D9') '0.0633'
$ws.Range('E9').Value = '  -0.78%  '
Set-TextValue $ws.Range('D10') '19.75'
$ws.Range('E10').Value = '  +0.79%  '
Set-TextValue $ws.Range('D11') '0.0787'
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('E12').Value = '  -0.35%  '
$ws.Range('D13').Value = '1.864.45'
$ws.Range('E13').Value = '  +0.25%  '
$ws.Range('D14').Value = '1.629.95'
$ws.Range('E14').Value = '  -0.33%  '
Set-TextValue $ws.Range('D15') '0.551'
$ws.Range('E15').Value = '  -1.34%  '
$ws.Range('E16').Value = '  -0.56%  '
Set-TextValue $ws.Range('D17') '63.20'
$ws.Range('E17').Value = '  +0.50%  '
$ws.Range('D18').Value = '26.043.75'
$ws.Range('E18').Value = '  +1.00%  '
$ws.Range('E19').Value = '  -0.91%  '
Set-TextValue $ws.Range('D20') '4.44'
$ws.Range('E20').Value = '  -0.36%  '
Set-TextValue $ws.Range('D21') '192.93'
$ws.Range('E21').Value = '  -0.77%  '
Set-TextValue $ws.Range('D22') '9.97'
$ws.Range('E22').Value = '  +0.04%  '
Set-TextValue $ws.Range('D23') '6.36'
$ws.Range('E23').Value = '  +1.27%  '
Set-TextValue $ws.Range('D24') '0.993'
$ws.Range('E24').Value = '  -0.95%  '
$ws.Range('E25').Value = '  -2.04%  '
Set-TextValue $ws.Range('D26') '141.60'
$ws.Range('E26').Value = '  -0.88%  '
Set-TextValue $ws.Range('D27') '0.124'
$ws.Range('E27').Value = '  +0.80%  '
Set-TextValue $ws.Range('D28') '6.88'
$ws.Range('E28').Value = '  -0.23%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('E30').Value = '  +0.41%  '
$ws.Range('E31').Value = '  -0.14%  '
$ws.Range('E32').Value = '  -0.09%  '
Set-TextValue $ws.Range('D33') '3.23'
$ws.Range('E33').Value = '  -0.75%  '
$ws.Range('E34').Value = '  +0.35%  '
Set-TextValue $ws.Range('D35') '2.39'
$ws.Range('E35').Value = '  -0.01%  '
Set-TextValue $ws.Range('D36') '0.906'
$ws.Range('E36').Value = '  +0.29%  '
$ws.Range('D37').Value = '1.141.95'
$ws.Range('E37').Value = '  +1.03%  '
$ws.Range('E38').Value = '  -0.42%  '
$ws.Range('E40').Value = '  -0.31%  '
Set-TextValue $ws.Range('D41') '0.993'
$ws.Range('E41').Value = '  -0.92%  '
$ws.Range('E42').Value = '  -0.56%  '
Set-TextValue $ws.Range('D43') '100.11'
$ws.Range('E43').Value = '  -0.12%  '
Set-TextValue $ws.Range('D44') '0.795'
$ws.Range('E44').Value = '  -1.43%  '
$ws.Range('D45').Value = '1.773.27'
$ws.Range('E45').Value = '  +0.22%  '
$ws.Range('D46').Value = '0.0₆0105'
$ws.Range('E46').Value = '  -3.40%  '
Set-TextValue $ws.Range('D47') '55.59'
$ws.Range('E47').Value = '  +0.72%  '
$ws.Range('E48').Value = '  +2.25%  '
$ws.Range('E49').Value = '  +4.37%  '
$ws.Range('E50').Value = '  -0.43%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D51') '0.0969'
$ws.Range('E51').Value = '  +2.70%  '
